# Update cryptos list — prices (column D) and volume/1h changes (column E)
# Values are stored as plain text in the sheet (coin prices such as
# "51.131.28" are not valid numbers and must stay textual; the volume
# percentages already carry padding spaces that keep them textual).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# row -> D(price) [$null = unchanged], E(volume 1h change)
$updates = @(
    @{ Row = 2;  D = "51.131.28";  E = "  -1.01%  " },
    @{ Row = 3;  D = "3.057.20";   E = "  +0.85%  " },
    @{ Row = 4;  D = $null;        E = "  +0.11%  " },
    @{ Row = 5;  D = "390.83";     E = "  +1.96%  " },
    @{ Row = 6;  D = "101.21";     E = "  -1.41%  " },
    @{ Row = 7;  D = "0.533";      E = "  -2.06%  " },
    @{ Row = 8;  D = $null;        E = "  +0.02%  " },
    @{ Row = 9;  D = "0.579";      E = "  -1.90%  " },
    @{ Row = 10; D = "36.62";      E = "  -0.22%  " },
    @{ Row = 11; D = $null;        E = "  +0.28%  " },
    @{ Row = 12; D = "0.0846";     E = "  -1.76%  " },
    @{ Row = 13; D = "3.545.54";   E = "  +1.03%  " },
    @{ Row = 14; D = "18.24";      E = "  -1.49%  " },
    @{ Row = 15; D = "7.64";       E = "  -1.15%  " },
    @{ Row = 16; D = "3.054.87";   E = "  +0.38%  " },
    @{ Row = 17; D = "0.990";      E = "  +1.91%  " },
    @{ Row = 18; D = "10.56";      E = "  -0.78%  " },
    @{ Row = 19; D = "51.138.67";  E = "  -1.03%  " },
    @{ Row = 20; D = $null;        E = "  +2.96%  " },
    @{ Row = 21; D = "12.21";      E = "  -2.09%  " },
    @{ Row = 22; D = "0.0₃0951";   E = "  -1.13%  " },
    @{ Row = 23; D = "69.59";      E = "  -0.55%  " },
    @{ Row = 24; D = "263.47";     E = "  -1.47%  " },
    @{ Row = 25; D = $null;        E = "  -1.11%  " },
    @{ Row = 26; D = "7.86";       E = "  -6.95%  " },
    @{ Row = 27; D = "26.69";      E = "  +1.69%  " },
    @{ Row = 28; D = "0.999";      E = "  -0.15%  " },
    @{ Row = 29; D = "7.11";       E = "  -4.97%  " },
    @{ Row = 30; D = "0.162";      E = "  -5.65%  " },
    @{ Row = 31; D = "0.104";      E = "  -3.34%  " },
    @{ Row = 32; D = "10.45";      E = "  +1.85%  " },
    @{ Row = 33; D = "0.0489";     E = "  +9.07%  " },
    @{ Row = 34; D = "35.59";      E = "  +4.30%  " },
    @{ Row = 35; D = "2.07";       E = "  -0.34%  " },
    @{ Row = 36; D = "49.93";      E = "  -1.28%  " },
    @{ Row = 37; D = $null;        E = "  +0.03%  " },
    @{ Row = 38; D = "3.34";       E = "  +0.84%  " },
    @{ Row = 39; D = "0.290";      E = "  -1.52%  " },
    @{ Row = 40; D = "129.38";     E = "  +1.14%  " },
    @{ Row = 41; D = "16.52";      E = "  -3.36%  " },
    @{ Row = 42; D = "1.82";       E = "  -2.49%  " },
    @{ Row = 43; D = "0.114";      E = "  -1.69%  " },
    @{ Row = 44; D = "3.76";       E = "  +1.50%  " },
    @{ Row = 45; D = "2.47";       E = "  -2.02%  " },
    @{ Row = 46; D = "21.63";      E = "  +0.08%  " },
    @{ Row = 47; D = $null;        E = "  +3.24%  " },
    @{ Row = 48; D = $null;        E = "  -0.24%  " },
    @{ Row = 49; D = "2.065.91";   E = "  +1.90%  " },
    @{ Row = 50; D = "0.0320";     E = "  +1.63%  " },
    @{ Row = 51; D = "0.882";      E = "  +10.80%  " }
)

foreach ($u in $updates) {
    $r = $u.Row

    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($r, 4)
        # Force the assigned value to stay text (prices like "51.131.28" or
        # "0.990" would otherwise be auto-coerced to numbers by Excel), then
        # restore the cell's original (default) style so no formatting
        # change is introduced.
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
        $dCell.Style = "Normal"
    }

    $eCell = $ws.Cells.Item($r, 5)
    $eCell.Value = $u.E
}
